$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Style = 'Normal'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.043.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.99%  '

# Row 3
$ws.Range('D3').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.013.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.76%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.41%  '

# Row 6
$ws.Range('E6').Value = '  -0.92%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('D8').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.29'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.68%  '

# Row 9
$ws.Range('E9').Value = '  -2.89%  '

# Row 10
$ws.Range('D10').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0777'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.13%  '

# Row 11
$ws.Range('E11').Value = '  -3.94%  '

# Row 12
$ws.Range('D12').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.307.64'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.87%  '

# Row 13
$ws.Range('D13').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.04'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.28%  '

# Row 14
$ws.Range('D14').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.27%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.736'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.39%  '

# Row 16
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.17'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.48%  '

# Row 17
$ws.Range('D17').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.012.55'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.96%  '

# Row 18
$ws.Range('D18').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.978.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.88%  '

# Row 19
$ws.Range('D19').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.44%  '

# Row 20
$ws.Range('D20').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.36%  '

# Row 21
$ws.Range('D21').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0811'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.61%  '

# Row 22
$ws.Range('D22').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.07%  '

# Row 23
$ws.Range('E23').Value = '  -0.02%  '

# Row 24
$ws.Range('E24').Value = '  +2.12%  '

# Row 25
$ws.Range('D25').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.71%  '

# Row 26
$ws.Range('D26').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.03%  '

# Row 27
$ws.Range('D27').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.47%  '

# Row 28
$ws.Range('D28').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.125'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.52%  '

# Row 29
$ws.Range('D29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.59'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.88%  '

# Row 30
$ws.Range('E30').Value = '  -7.35%  '

# Row 31
$ws.Range('D31').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.117'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.04%  '

# Row 32
$ws.Range('D32').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.75%  '

# Row 33
$ws.Range('E33').Value = '  -1.99%  '

# Row 34
$ws.Range('E34').Value = '  -2.34%  '

# Row 35
$ws.Range('D35').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.33'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.74%  '

# Row 36
$ws.Range('E36').Value = '  +2.13%  '

# Row 37
$ws.Range('D37').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.24%  '

# Row 38
$ws.Range('D38').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.10%  '

# Row 39
$ws.Range('E39').Value = '  -1.39%  '

# Row 40
$ws.Range('D40').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.462.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.00%  '

# Row 41
$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.31'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +19.99%  '

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0212'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.28%  '

# Row 43
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '94.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.76%  '

# Row 44
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0909'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.60%  '

# Row 45
$ws.Range('E45').Value = '  -4.09%  '

# Row 46
$ws.Range('E46').Value = '  -2.05%  '

# Row 47
$ws.Range('D47').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.32%  '

# Row 48
$ws.Range('D48').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.997'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.99%  '

# Row 49
$ws.Range('D49').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.07'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.31%  '

# Row 50
$ws.Range('E50').Value = '  -1.04%  '

# Row 51
$ws.Range('D51').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.195.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.84%  '
